# Updated symbol list on Sat Jan 14 08:43:55 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows that moved since the last snapshot. Values are written as
# literal text (NumberFormat "@") so that strings like "8.70%" or
# "2,567.65%" are preserved exactly rather than being re-interpreted by
# Excel as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '312.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '8.70%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '32.47'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '9.72%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.337'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.22%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07678'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '14.72%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.879'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '7.48%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.711'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '8.98%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.613'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '18.42%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9192'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.71%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01724'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2,567.65%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1724'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '8.47%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07573'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '12.24%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08250'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '7.24%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03027'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.16%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09885'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '10.18%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001526'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.90%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04560'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.54%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006145'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.75%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.477'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.01%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.244'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.07%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.13%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1307'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.13%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.246'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.33%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.78%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.72%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004507'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.50%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001298'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '8.25%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001737'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '7.47%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04650'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '9.06%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007216'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '6.70%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1372'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '10.79%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002256'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.26%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01461'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.93%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006211'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '9.57%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.893'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-3.82%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01297'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.65%'
